# Motilal portfolio change engine — insert an "Industry" column after
# "Stock Name" (new column C), pushing "Mutual Fund" ... "QoQ" one column
# to the right (D..J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns C:I -> D:J by inserting a fresh column at C.
$ws.Columns("C").Insert()

# New header cell for the inserted column.
$ws.Range("C1").Value = "Industry"

# Per-row Industry values for the newly inserted column C (rows 2-21).
$industry = @{
    2  = "IT - Services"
    3  = "IT - Software"
    4  = "IT - Software"
    5  = "IT - Services"
    6  = "IT - Software"
    7  = "Metals & Minerals Trading"
    8  = "IT - Software"
    9  = "IT - Software"
    10 = "IT - Software"
    11 = "IT - Software"
    12 = "IT - Software"
    13 = "IT - Software"
    14 = "Banks"
    15 = "Entertainment"
    16 = "Retailing"
    17 = "Telecom - Services"
    18 = "IT - Services"
    19 = "IT - Services"
    20 = "IT - Software"
    21 = "Commercial Services & Supplies"
}

foreach ($row in $industry.Keys) {
    $ws.Cells.Item($row, 3).Value = $industry[$row]
}
